$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$remark = "दिनांक 02.12.2025 रोजी रविवार असल्याने जमा झालेली रक्कम ही सोमवार दिनांक 03.12.2025 रोजी बँकेमध्ये भरणा करण्यात आली."

# Update Remark column (F) for rows 7-10 to include the date
$ws.Range("F7").Value = $remark
$ws.Range("F8").Value = $remark
$ws.Range("F9").Value = $remark
$ws.Range("F10").Value = $remark

# Update Mode column (C) for rows 8-11
$ws.Range("C8").Value = "Cash"
$ws.Range("C9").Value = "Cheque"
$ws.Range("C10").Value = "NEFT"
$ws.Range("C11").Value = "Total"
